# Natmi following Dr Hou advice
# Update Egf-Erbb3 LR-pair sheet: add a third cluster ("ECs") to the
# existing two ("FAPs", "sCs"), expanding the sending/target cluster
# combinations from a 2x2 matrix (rows 2-5) to a 3x3 matrix (rows 2-10),
# and refresh every numeric metric column (E:T) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Egf"
$ws.Range("C2").Value = "Erbb3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.050561
$ws.Range("H2").Value = 0.151683
$ws.Range("I2").Value = 0.1845256053410153
$ws.Range("J2").Value = 0.1845256053410153
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.058393
$ws.Range("N2").Value = 0.175179
$ws.Range("O2").Value = 0.01085373024912483
$ws.Range("P2").Value = 0.01085373024912483
$ws.Range("Q2").Value = 0.002952408473
$ws.Range("R2").Value = 0.026571676257
$ws.Range("S2").Value = 0.002002791144427848
$ws.Range("T2").Value = 0.002002791144427848
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Egf"
$ws.Range("C3").Value = "Erbb3"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.050561
$ws.Range("H3").Value = 0.151683
$ws.Range("I3").Value = 0.1845256053410153
$ws.Range("J3").Value = 0.1845256053410153
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4773683333333333
$ws.Range("N3").Value = 1.432105
$ws.Range("O3").Value = 0.08873027793527143
$ws.Range("P3").Value = 0.08873027793527143
$ws.Range("Q3").Value = 0.02413622030166667
$ws.Range("R3").Value = 0.217225982715
$ws.Range("S3").Value = 0.01637300824808249
$ws.Range("T3").Value = 0.01637300824808249
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Egf"
$ws.Range("C4").Value = "Erbb3"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.050561
$ws.Range("H4").Value = 0.151683
$ws.Range("I4").Value = 0.1845256053410153
$ws.Range("J4").Value = 0.1845256053410153
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.844232333333333
$ws.Range("N4").Value = 14.532697
$ws.Range("O4").Value = 0.9004159918156037
$ws.Range("P4").Value = 0.9004159918156038
$ws.Range("Q4").Value = 0.2449292310056667
$ws.Range("R4").Value = 2.204363079051
$ws.Range("S4").Value = 0.1661498059485049
$ws.Range("T4").Value = 0.166149805948505
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Egf"
$ws.Range("C5").Value = "Erbb3"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1101723333333333
$ws.Range("H5").Value = 0.3305169999999999
$ws.Range("I5").Value = 0.4020809813920896
$ws.Range("J5").Value = 0.4020809813920896
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.058393
$ws.Range("N5").Value = 0.175179
$ws.Range("O5").Value = 0.01085373024912483
$ws.Range("P5").Value = 0.01085373024912483
$ws.Range("Q5").Value = 0.006433293060333333
$ws.Range("R5").Value = 0.05789963754299999
$ws.Range("S5").Value = 0.004364078510333121
$ws.Range("T5").Value = 0.004364078510333121
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Egf"
$ws.Range("C6").Value = "Erbb3"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1101723333333333
$ws.Range("H6").Value = 0.3305169999999999
$ws.Range("I6").Value = 0.4020809813920896
$ws.Range("J6").Value = 0.4020809813920896
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4773683333333333
$ws.Range("N6").Value = 1.432105
$ws.Range("O6").Value = 0.08873027793527143
$ws.Range("P6").Value = 0.08873027793527143
$ws.Range("Q6").Value = 0.05259278314277777
$ws.Range("R6").Value = 0.4733350482849999
$ws.Range("S6").Value = 0.03567675723140681
$ws.Range("T6").Value = 0.03567675723140681
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Egf"
$ws.Range("C7").Value = "Erbb3"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1101723333333333
$ws.Range("H7").Value = 0.3305169999999999
$ws.Range("I7").Value = 0.4020809813920896
$ws.Range("J7").Value = 0.4020809813920896
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.844232333333333
$ws.Range("N7").Value = 14.532697
$ws.Range("O7").Value = 0.9004159918156037
$ws.Range("P7").Value = 0.9004159918156038
$ws.Range("Q7").Value = 0.533700379372111
$ws.Range("R7").Value = 4.803303414348999
$ws.Range("S7").Value = 0.3620401456503497
$ws.Range("T7").Value = 0.3620401456503498
$ws.Range("A8").Value = "ECs"
$ws.Range("B8").Value = "Egf"
$ws.Range("C8").Value = "Erbb3"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.113272
$ws.Range("H8").Value = 0.339816
$ws.Range("I8").Value = 0.413393413266895
$ws.Range("J8").Value = 0.413393413266895
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.058393
$ws.Range("N8").Value = 0.175179
$ws.Range("O8").Value = 0.01085373024912483
$ws.Range("P8").Value = 0.01085373024912483
$ws.Range("Q8").Value = 0.006614291896
$ws.Range("R8").Value = 0.059528627064
$ws.Range("S8").Value = 0.00448686059436386
$ws.Range("T8").Value = 0.00448686059436386
$ws.Range("A9").Value = "ECs"
$ws.Range("B9").Value = "Egf"
$ws.Range("C9").Value = "Erbb3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.113272
$ws.Range("H9").Value = 0.339816
$ws.Range("I9").Value = 0.413393413266895
$ws.Range("J9").Value = 0.413393413266895
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4773683333333333
$ws.Range("N9").Value = 1.432105
$ws.Range("O9").Value = 0.08873027793527143
$ws.Range("P9").Value = 0.08873027793527143
$ws.Range("Q9").Value = 0.05407246585333333
$ws.Range("R9").Value = 0.48665219268
$ws.Range("S9").Value = 0.03668051245578212
$ws.Range("T9").Value = 0.03668051245578212
$ws.Range("A10").Value = "ECs"
$ws.Range("B10").Value = "Egf"
$ws.Range("C10").Value = "Erbb3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.113272
$ws.Range("H10").Value = 0.339816
$ws.Range("I10").Value = 0.413393413266895
$ws.Range("J10").Value = 0.413393413266895
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.844232333333333
$ws.Range("N10").Value = 14.532697
$ws.Range("O10").Value = 0.9004159918156037
$ws.Range("P10").Value = 0.9004159918156038
$ws.Range("Q10").Value = 0.5487158848613333
$ws.Range("R10").Value = 4.938442963752
$ws.Range("S10").Value = 0.372226040216749
$ws.Range("T10").Value = 0.3722260402167491
